$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'2019"
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A3").Value = "'2019"
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
